$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing row 7 (Q5) values
$ws.Range("B7").Value = -0.4243341782850507
$ws.Range("C7").Value = 1.33319951928904
$ws.Range("D7").Value = 2.878316400685
$ws.Range("E7").Value = 1.696560167127886
$ws.Range("F7").Value = 1.74227983373859
$ws.Range("G7").Value = 9

# Update existing row 8 (Q6) values
$ws.Range("B8").Value = -0.8548487235232874
$ws.Range("C8").Value = 1.0621450549565
$ws.Range("D8").Value = 1.568653354045059
$ws.Range("E8").Value = 1.252458923096905
$ws.Range("F8").Value = 1.002728486043355
$ws.Range("G8").Value = 6

# Update existing row 9 (Q7) values, and add new F9 value
$ws.Range("B9").Value = -0.4816660954554475
$ws.Range("C9").Value = 1.068906077927655
$ws.Range("D9").Value = 1.932725749686544
$ws.Range("E9").Value = 1.390225071593281
$ws.Range("F9").Value = 1.597211721489318
$ws.Range("G9").Value = 3

# Add new row 10 (Q8), copying the format of row 9's label cell for A10
$ws.Range("A10").Value = "Q8"
$ws.Range("A9").Copy()
$ws.Range("A10").PasteSpecial(-4122)
$ws.Range("B10").Value = -0.1119550751434417
$ws.Range("C10").Value = 0.1119550751434417
$ws.Range("D10").Value = 0.01253393885037368
$ws.Range("E10").Value = 0.1119550751434417
$ws.Range("G10").Value = 1
